# Horarios actualizados Linea 141 - 10:12:35 scrape run
# Updates LP1912, LP1912-215 and 6203-6173 sheets with the newly scraped
# arrival rows (re-sorted by Hora_Llegada) and refreshes the summary header.

$wb = $excel.ActiveWorkbook

function Set-ScheduleRow($ws, $r, $horaScrap, $horaLlegada, $linea, $minutos, $parada) {
    $ws.Cells.Item($r, 1).Value = $horaScrap
    $ws.Cells.Item($r, 2).Value = $horaLlegada
    $ws.Cells.Item($r, 3).Value = $linea
    $ws.Cells.Item($r, 4).Value = $minutos
    $ws.Cells.Item($r, 5).Value = $parada
}

$ws = $wb.Worksheets.Item('LP1912')
$ws.Range('A2').Value = 'Última actualización: 10:12:35'
$ws.Range('A3').Value = 'Total filas: 157'

$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range('A2').Value = 'Última actualización: 10:12:35'
$ws.Range('A3').Value = 'Total filas: 20'

$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range('A2').Value = 'Última actualización: 10:12:35'
$ws.Range('A3').Value = 'Total filas: 30'

$ws = $wb.Worksheets.Item('LP1912')
Set-ScheduleRow $ws 133 '10:12:35' '10:33' '10_OLMOS' 21 'LP1912'
Set-ScheduleRow $ws 134 '10:12:35' '10:34' '16_SANTA ANA' 22 'LP1912'
Set-ScheduleRow $ws 135 '10:12:35' '10:35' '23_HERNANDEZ' 23 'LP1912'
Set-ScheduleRow $ws 136 '09:25:56' '10:41' '17_ROMERO' 76 'LP1912'
Set-ScheduleRow $ws 137 '08:48:09' '10:42' '17_ROMERO' 114 'LP1912'
Set-ScheduleRow $ws 138 '08:55:19' '10:43' '14_ABASTO' 108 'LP1912'
Set-ScheduleRow $ws 139 '08:48:09' '10:44' '14_ABASTO' 116 'LP1912'
Set-ScheduleRow $ws 140 '10:12:35' '10:52' '15_ABASTO' 40 'LP1912'
Set-ScheduleRow $ws 141 '09:25:56' '10:53' '27_EL RETIRO' 88 'LP1912'
Set-ScheduleRow $ws 142 '10:12:35' '10:53' '10_OLMOS' 41 'LP1912'
Set-ScheduleRow $ws 143 '10:12:35' '10:56' '27_EL RETIRO' 44 'LP1912'
Set-ScheduleRow $ws 144 '09:25:56' '11:01' '215C_EL PATO' 96 'LP1912'
Set-ScheduleRow $ws 145 '10:12:35' '11:02' '215C_EL PATO' 50 'LP1912'
Set-ScheduleRow $ws 146 '10:12:35' '11:03' '11_ETCHEVERRY' 51 'LP1912'
Set-ScheduleRow $ws 147 '10:12:35' '11:04' '23_HERNANDEZ' 52 'LP1912'
Set-ScheduleRow $ws 148 '09:25:56' '11:10' '16_P MOR-167 Y 521' 105 'LP1912'
Set-ScheduleRow $ws 149 '10:12:35' '11:12' '15_ABASTO' 60 'LP1912'
Set-ScheduleRow $ws 150 '09:25:56' '11:19' '86_EST CHICA-ESC AGRARIA' 114 'LP1912'
Set-ScheduleRow $ws 151 '09:25:56' '11:20' '26_HERNANDEZ' 115 'LP1912'
Set-ScheduleRow $ws 152 '10:12:35' '11:21' '26_HERNANDEZ' 69 'LP1912'
Set-ScheduleRow $ws 153 '10:12:35' '11:27' '225_C ROCA-H SUR' 75 'LP1912'
Set-ScheduleRow $ws 154 '10:12:35' '11:32' '81_EL PELIGRO' 80 'LP1912'
Set-ScheduleRow $ws 155 '10:12:35' '11:42' '17_ROMERO' 90 'LP1912'
Set-ScheduleRow $ws 156 '10:12:35' '11:51' '215B_EL PATO' 99 'LP1912'
Set-ScheduleRow $ws 157 '10:12:35' '11:56' '10_OLMOS' 104 'LP1912'
Set-ScheduleRow $ws 158 '10:12:35' '11:58' '16_P MOR-167 Y 521' 106 'LP1912'
Set-ScheduleRow $ws 159 '10:12:35' '11:59' '225_GOMEZ' 107 'LP1912'
Set-ScheduleRow $ws 160 '10:12:35' '12:04' '84_COLONIA URQUIZA-ESC 49' 112 'LP1912'
Set-ScheduleRow $ws 161 '10:12:35' '12:06' '16_P MOR-SANTA ANA' 114 'LP1912'
Set-ScheduleRow $ws 162 '10:12:35' '12:06' '14_ABASTO' 114 'LP1912'

$ws = $wb.Worksheets.Item('LP1912-215')
Set-ScheduleRow $ws 24 '10:12:35' '11:02' '215C_EL PATO' 50 'LP1912'
Set-ScheduleRow $ws 25 '10:12:35' '11:51' '215B_EL PATO' 99 'LP1912'

$ws = $wb.Worksheets.Item('6203-6173')
Set-ScheduleRow $ws 35 '10:12:35' '12:04' '215A_LA PLATA' 112 'L6173'
